$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update InsCardNo value in X2 (shared string change)
$ws.Range("X2").Value = "DN4127389127641"

# Update row 2: A2 (id) and E2 (some code)
$ws.Range("A2").Value = 199
$ws.Range("E2").Value = 46200020983

# Update row 3: A3 (id) and E3 (some code)
$ws.Range("A3").Value = 200
$ws.Range("E3").Value = 46200020984

# Update the active selection to J11
$ws.Range("J11").Select()
